# bill invoice -> credit purchase and product import bugfix
#
# The sample_items product-import template gets a new leading "id" column,
# and the misspelled "intial_stock" header is fixed to "initial_stock".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts every existing column
# (data, styles, and column-width formatting) one place to the right.
$ws.Columns("A:A").Insert()

# New leading column header.
$ws.Range("A1").Value = "id"

# Fix the "intial_stock" -> "initial_stock" typo (this column was G, now H
# after the insert above).
$ws.Range("H1").Value = "initial_stock"

# Match the saved selection/active cell shown in the target workbook.
$ws.Range("H2").Select() | Out-Null
